# Set_9 (all.gp = TRUE) params
# Fill in the "Corrected_analysis_*" columns (L:O) for rows 14, 16 and 17
# on the "Sims and main analysis" sheet, and update the active
# cell/selection to match where the author last left the cursor.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sims and main analysis")

# --- Row 14 (Set_13) -------------------------------------------------
$ws.Range("L14").Value = 45800
$ws.Range("L14").NumberFormat = "d-mmm-yy"
$ws.Range("M14").Value = "48782249 (Eddie)"
$ws.Range("N14").Value = 45802
$ws.Range("N14").NumberFormat = "d-mmm-yy"
$ws.Range("O14").Value = "NA"

# --- Row 16 (Set_N1 (157-168)) ----------------------------------------
$ws.Range("O16").Value = "NA"

# --- Row 17 (Set_N1 (169-176)) ----------------------------------------
$ws.Range("O17").Value = "NA"

# --- Selection / active cell ------------------------------------------
# The frozen-pane split (column A frozen, xSplit = 1) is unchanged;
# only the active cell moved.
$ws.Range("O17").Select()
